$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row -> new "Data" (column E) text value, for years 1950-2016 (rows 2-68).
# Rows 2-60 already exist (years 1950-2008); rows 61-68 are brand new
# (years 2009-2016) and need the full row (A-E) populated.
$yearData = @(
    @{ Row = 2; Val = "631" },
    @{ Row = 3; Val = "657" },
    @{ Row = 4; Val = "666" },
    @{ Row = 5; Val = "682" },
    @{ Row = 6; Val = "708" },
    @{ Row = 7; Val = "719" },
    @{ Row = 8; Val = "735" },
    @{ Row = 9; Val = "749" },
    @{ Row = 10; Val = "764" },
    @{ Row = 11; Val = "773" },
    @{ Row = 12; Val = "813" },
    @{ Row = 13; Val = "813" },
    @{ Row = 14; Val = "918" },
    @{ Row = 15; Val = "995" },
    @{ Row = 16; Val = "1054" },
    @{ Row = 17; Val = "1054" },
    @{ Row = 18; Val = "1025" },
    @{ Row = 19; Val = "1113" },
    @{ Row = 20; Val = "1084" },
    @{ Row = 21; Val = "1074" },
    @{ Row = 22; Val = "1073" },
    @{ Row = 23; Val = "952" },
    @{ Row = 24; Val = "1108" },
    @{ Row = 25; Val = "1364" },
    @{ Row = 26; Val = "1414" },
    @{ Row = 27; Val = "1270" },
    @{ Row = 28; Val = "1447" },
    @{ Row = 29; Val = "1658" },
    @{ Row = 30; Val = "1934" },
    @{ Row = 31; Val = "1678" },
    @{ Row = 32; Val = "1589" },
    @{ Row = 33; Val = "1564" },
    @{ Row = 34; Val = "1576" },
    @{ Row = 35; Val = "1561" },
    @{ Row = 36; Val = "1589" },
    @{ Row = 37; Val = "1600" },
    @{ Row = 38; Val = "1645" },
    @{ Row = 39; Val = "1621" },
    @{ Row = 40; Val = "1736" },
    @{ Row = 41; Val = "1804" },
    @{ Row = 42; Val = "1882" },
    @{ Row = 43; Val = "1853.54297901117" },
    @{ Row = 44; Val = "1914.62382045181" },
    @{ Row = 45; Val = "1903.67622257841" },
    @{ Row = 46; Val = "1924.09465245235" },
    @{ Row = 47; Val = "1888.69902048355" },
    @{ Row = 48; Val = "1931.6607745291" },
    @{ Row = 49; Val = "1940.31058051924" },
    @{ Row = 50; Val = "1922.67938128057" },
    @{ Row = 51; Val = "1923.90323254743" },
    @{ Row = 52; Val = "1967.43327436145" },
    @{ Row = 53; Val = "1997.35760955142" },
    @{ Row = 54; Val = "1959.02789399239" },
    @{ Row = 55; Val = "2003.08991926606" },
    @{ Row = 56; Val = "2004.06017812124" },
    @{ Row = 57; Val = "2013.80775709797" },
    @{ Row = 58; Val = "2056.32966356963" },
    @{ Row = 59; Val = "2106.84014026248" },
    @{ Row = 60; Val = "2177.85735056552" },
    @{ Row = 61; Val = "2199.39485848654" },
    @{ Row = 62; Val = "2313.88212836641" },
    @{ Row = 63; Val = "2346" },
    @{ Row = 64; Val = "2456" },
    @{ Row = 65; Val = "2559" },
    @{ Row = 66; Val = "2644" },
    @{ Row = 67; Val = "2710" },
    @{ Row = 68; Val = "2767" }
)

foreach ($item in $yearData) {
    $row = $item.Row
    $year = 1950 + ($row - 2)

    if ($row -ge 61) {
        # New row: Country Code, Country Name, Indicator, Year
        $ws.Cells.Item($row, 1).Value = 426
        $ws.Cells.Item($row, 2).Value = "Lesotho"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Column E ("Data") must stay a text value (it looks numeric but the
    # source data stores it as shared-string text), so force the Text
    # number format before assigning, then clear the format again so the
    # cell keeps the workbook's default (General) style.
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Val
    $cell.ClearFormats()
}
